# forBot_StudentRequestsMatrix.xlsx — "Fix for students who request no one"
#
# 1) Several faculty column headers were abbreviated to last-name-only;
#    restore full names.
# 2) A student ("Shen, Jiayi") who had requested nobody is removed entirely.
# 3) A new faculty column ("John Lowengrub") is inserted (between "Chang Liu"
#    and "Ilhem Messaoudi") and a few students' requests now include him.
# 4) Two students ("Du, Mingyu" and "Xiang, Yankai (Mark)") who previously had
#    a single (apparently erroneous) request now request no one — their single
#    "Jun Allard" tick is cleared.
# 5) Two other students ("Hwang, Ahyeon" and "Wang, Harold") lose an erroneous
#    "Lander" (Arthur Lander) tick.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Fix abbreviated faculty header names (while columns are still in
#        their original positions) ---------------------------------------
$ws.Range("E1").Value = "Lee Bardwell"
$ws.Range("G1").Value = "James Brody"
$ws.Range("N1").Value = "Arthur Lander"
$ws.Range("T1").Value = "Ali Mortazavi"
$ws.Range("W1").Value = "Jose Ranz"
$ws.Range("X1").Value = "Thomas Schilling"
$ws.Range("Z1").Value = "Albert Siryaporn"
$ws.Range("AB1").Value = "Vivek Swarup"
$ws.Range("AE1").Value = "Katrine Whiteson"
$ws.Range("AF1").Value = "Dominik Wodarz"
$ws.Range("AG1").Value = "Zeba Wunderlich"

# --- 2) Data corrections for existing students (still original row/col
#        positions) --------------------------------------------------------
# Du, Mingyu (row 7) and Xiang, Yankai (Mark) (row 26) now request no one.
$ws.Range("B7").Value = 0
$ws.Range("B26").Value = 0

# Hwang, Ahyeon (row 14) and Wang, Harold (row 24) no longer request
# Arthur Lander (col N).
$ws.Range("N14").Value = 0
$ws.Range("N24").Value = 0

# --- 3) Remove the student who requested no one: "Shen, Jiayi" (row 21) ---
$ws.Rows(21).Delete()

# --- 4) Insert the new faculty column "John Lowengrub" between "Chang Liu"
#        (Q) and "Ilhem Messaoudi" (R) --------------------------------------
$ws.Columns("R").Insert()
$ws.Range("R1").Value = "John Lowengrub"

# Default the whole new column's data rows to 0, then flip on the students
# who request him.
$ws.Range("R2:R27").Value = 0
$ws.Range("R19").Value = 1   # Narain, Vedang
$ws.Range("R20").Value = 1   # Nguyen, Nguyen
$ws.Range("R22").Value = 1   # Sousa, Rachel
$ws.Range("R24").Value = 1   # Ward, Erica
